$d = $word.ActiveDocument
$find = "ensure the long-term success of their IT infrastructure"
$replace = "ensure the long-term success of their IT infrastructure."
$d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
